$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update C2 last-updated date for Daily measurements
$ws.Range("C2").Value = 20221012

# Add new row 6: Feeding protocol entry
$ws.Range("A6").Value = "Feeding"
$ws.Range("B6").Value = "https://github.com/urol-e5/protocols/blob/master/2022-10-17-Moorea-E5-Feeding-Protocol.md"
$ws.Range("C6").Value = 20221017

# Update selection to match the final state (B11)
$ws.Range("B11").Select()
